# Persona Description.docx
# Move the "_GoBack" bookmark out of the paragraph right after "Personal
# Profile:" and append a new "Scenario of use: " section (with two blank
# paragraphs above it, Overskrift2 style, and the bookmark) at the very
# end of the document body, after the existing "Business objectives:"
# paragraph.

$d = $word.ActiveDocument

# 1) Remove the bookmark from its old location (empty paragraph right
#    after the "Personal Profile:" heading). The empty paragraph itself
#    is left untouched.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Build the new content: two blank paragraphs followed by a new
#    "Overskrift2"-styled heading paragraph "Scenario of use: " which
#    carries the relocated "_GoBack" bookmark (collapsed, right after the
#    text, as in the original document).
$newContentXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Overskrift2"/>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">Scenario of use: </w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

# 3) Insert the new content at the very end of the document body (after
#    the last paragraph, "Business objectives:").
$endOfDoc = $d.Content
$endOfDoc.Collapse(0)
$endOfDoc.InsertXML($newContentXml)
